# "update2 of excel skills"
# Insert a new "R-library" skill row ("dataframes treatment") right above the
# existing "Lubridate" row (row 11) on the Skills sheet, pushing the rest of
# the R-library block (and the trailing blank rows) down by one row, and
# update the sheet view/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row at position 11; Excel shifts rows 11..25 down to
# 12..26 automatically, carrying their content/styles/borders with them.
$ws.Rows.Item(11).Insert()

# The freshly inserted row 11 has no formatting yet (cells are blank).
# Copy the row-formatting (borders/fills/fonts/number formats) from the row
# right below it (row 12, which now holds what used to be row 11 - the
# "Lubridate" entry) so the new row matches the rest of the table.
$ws.Range("A12:F12").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column B in this table highlights a couple of entries with a yellow fill
# (style carried over by the paste above). The new row should NOT be
# highlighted, so copy the plain (non-highlighted) formatting from A11
# onto B11 instead.
$ws.Range("A11").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fill in the new row's values.
$ws.Range("A11").Value = "R-library"
$ws.Range("B11").Value = "dataframes treatment"
$ws.Range("F11").Value = "x"

# Update the view: select E12 as the active cell (and let the top-left cell
# reset to the default/top of the sheet).
$null = $ws.Activate()
$null = $ws.Range("E12").Select()
